$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Remove the final three slides:
#      8  "Capturing and versioning Data and Analysis Meta Data"
#      9  "Modification to Project Structure"
#      10 "Exercise"
#    (delete from the highest index down so earlier indices stay valid)
# ---------------------------------------------------------------------------
$p.Slides.Item(10).Delete()
$p.Slides.Item(9).Delete()
$p.Slides.Item(8).Delete()

# ---------------------------------------------------------------------------
# 2. On the "Project organization" content slide (still slide 5), consolidate
#    the runs that make up the first two bullets into single runs.
# ---------------------------------------------------------------------------
$s   = $p.Slides.Item(5)
$shp = $s.Shapes.Item(2)
$tr  = $shp.TextFrame.TextRange

# Paragraph 1: "store all " + "of the files relevant to one " + "project under " +
#              "a common root " + "directory"  ->  single run
$para1 = $tr.Paragraphs(1, 1)
$c1 = $tr.Characters($para1.Start, $para1.Length - 1)
$c1.Text = "store all of the files relevant to one project under a common root directory"

# Paragraph 2: "logical " + "top-" + "level organization" -> single run
$tr2 = $shp.TextFrame.TextRange
$para2 = $tr2.Paragraphs(2, 1)
$c2 = $tr2.Characters($para2.Start, $para2.Length - 1)
$c2.Text = "logical top-level organization"

# ---------------------------------------------------------------------------
# 3. Drop the stray trailing endParaRPr after "logical tertiary organization"
#    by removing that paragraph and retyping it right after the previous one.
# ---------------------------------------------------------------------------
$tr3 = $shp.TextFrame.TextRange
$para9 = $tr3.Paragraphs(9, 1)
$c3 = $tr3.Characters($para9.Start, $para9.Length - 1)
$c3.Text = ""

$tr4 = $shp.TextFrame.TextRange
$para8 = $tr4.Paragraphs(8, 1)
$para8.InsertAfter("`rlogical tertiary organization")
